$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) holds a date serial value that was updated from
# 45793 to 45794 for every data row (rows 2 through 43).
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45793) {
        $cell.Value = 45794
    }
}
